# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet so the existing "Late" / "Outstanding" data block shifts one
# column to the right (N->O, O(blank)->P, P->Q), matching the new
# "Variable Instalments" layout, then restore the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column at position N (14). Everything from N onward
# (including the header cells and the per-row values) shifts right by one.
$ws.Columns.Item(14).Insert()

# The newly inserted column inherits the look of its left neighbour (M);
# give it the same width used by the rest of the schedule's value columns.
$ws.Columns.Item(14).ColumnWidth = 10.25

# Excel leaves the current selection on the sheet that was active when the
# file was last saved - match the author's recorded selection.
$ws.Range("R7").Select()
